# Add a new "BB" column to the right of the existing "BA" column.
# BB1 gets the next period date-serial; BB3..BB21 mostly carry forward the
# BA value for that row (latest forecast value), except rows 19-21 which
# get refreshed/updated forecast figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NewColumnCell($targetRef, $sourceRef, $value) {
    # Copy the source cell's formatting (number format/font/border/etc.)
    # onto the target cell so it reuses the existing style, then set the
    # numeric value explicitly.
    $ws.Range($sourceRef).Copy() | Out-Null
    $ws.Range($targetRef).PasteSpecial(-4122) | Out-Null
    $ws.Range($targetRef).Value = $value
}

# Row 1: header date serial
Set-NewColumnCell "BB1" "BA1" 45986

# Rows 3-18: new column simply repeats the last (BA) forecast value
Set-NewColumnCell "BB3"  "BA3"  -14.1223525803845
Set-NewColumnCell "BB4"  "BA4"  7.235571181590705
Set-NewColumnCell "BB5"  "BA5"  9.488155060723313
Set-NewColumnCell "BB6"  "BA6"  4.445724792911898
Set-NewColumnCell "BB7"  "BA7"  1.101453765199745
Set-NewColumnCell "BB8"  "BA8"  3.917168917088798
Set-NewColumnCell "BB9"  "BA9"  4.449305425651406
Set-NewColumnCell "BB10" "BA10" 1.92981802270098
Set-NewColumnCell "BB11" "BA11" 4.693514706708668
Set-NewColumnCell "BB12" "BA12" 4.873158811425493
Set-NewColumnCell "BB13" "BA13" 0.862099696137153
Set-NewColumnCell "BB14" "BA14" -8.794825048137467
Set-NewColumnCell "BB15" "BA15" 6.095721945083143
Set-NewColumnCell "BB16" "BA16" 5.588511616267167
Set-NewColumnCell "BB17" "BA17" -0.008466500317649839
Set-NewColumnCell "BB18" "BA18" -1.48753958890171

# Rows 19-21: new column carries an updated forecast value
Set-NewColumnCell "BB19" "BA19" -2.451276118722334
Set-NewColumnCell "BB20" "BA20" 1.795477855501626
Set-NewColumnCell "BB21" "BA21" 1.878198916198426
